$wb = $excel.ActiveWorkbook

# --- "chambers" sheet: rename headers (b_/s_ -> chb_/chs_) and fill in zeros for row 4 ---
$wsChambers = $wb.Worksheets.Item("chambers")
$wsChambers.Range("D1").Value = "chb_irr1"
$wsChambers.Range("E1").Value = "chb_con2"
$wsChambers.Range("F1").Value = "chs_irr3"
$wsChambers.Range("G1").Value = "chs_con4"
$wsChambers.Range("H1").Value = "chs_irr5"
$wsChambers.Range("I1").Value = "chs_con6"

$wsChambers.Range("D4").Value = 0
$wsChambers.Range("E4").Value = 0
$wsChambers.Range("F4").Value = 0
$wsChambers.Range("G4").Value = 0

# --- "flow" sheet: fill in previously-blank data row 3 ---
$wsFlow = $wb.Worksheets.Item("flow")
$wsFlow.Range("A3").Value = 44573
$wsFlow.Range("B3").Value = 0
$wsFlow.Range("C3").Value = 0
$wsFlow.Range("D3").Value = 1
$wsFlow.Range("E3").Value = 2
$wsFlow.Range("F3").Value = 0
$wsFlow.Range("G3").Value = 0

# --- Update the selections / active cells on each sheet ---
$wsPar = $wb.Worksheets.Item("par")
$wsPar.Range("A3").Select()

$wsFlow.Range("G4").Select()

# "chambers" becomes the active sheet/tab, with H4 as the selected cell
$wsChambers.Range("H4").Select()
$wsChambers.Activate()
